$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# Keep the replacement stats (runs/balls/4s/6s/SR) stored as text, same as
# the rest of the sheet, rather than letting Excel auto-detect them as numbers
$ws.Range("G2:K2").NumberFormat = "@"

# Row 2 becomes the data that used to live in row 4 (the "Mumbai" match) -
# all the other historical rows for this player are being dropped
$ws.Range("A2").Value = " Oct 16 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Eoin Morgan${nbsp}(c)"
$ws.Range("G2").Value = "39"
$ws.Range("H2").Value = "29"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "2"
$ws.Range("K2").Value = "134.48"

# Remove the now-stale rows 3 through 8 so only the header + the single
# remaining match row are left (sheet dimension collapses to A1:K2)
$ws.Range("A3:K8").Delete()
